# aggiornamento fino a 02/05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates as Excel serials, matching the existing A-column style)
$data = @(
    @(44313, 1, 7, 80.49678012879485),
    @(44314, 0, 7, 80.49678012879485),
    @(44315, 3, 9, 103.4958601655934),
    @(44316, 0, 8, 91.99632014719411),
    @(44317, 0, 6, 68.99724011039559),
    @(44318, 0, 4, 45.99816007359706)
)

$startRow = 239
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Copy the style of the previous row's cells so the new rows render the
    # same way (date format + border on column A, plain on B/C/D).
    $ws.Cells.Item($r - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
